$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55, column A currently holds the text "09876543" (leading zero).
# The edit corrects it to the plain number 9876543.
$ws.Range("A55").Value = 9876543

# Append the new payment record as row 56.
# Column A is a phone-style id with a leading zero, so it must stay text.
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "09876543"
$ws.Range("A56").Style = "Normal"
$ws.Range("B56").Value = ""
$ws.Range("C56").Value = "Cash"
$ws.Range("D56").Value = "2025-08-18T18:04:40"
$ws.Range("E56").Value = 120
$ws.Range("F56").Value = ""
$ws.Range("G56").Value = 120
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
